# Fix a student-reported typo on slide 28 ("lec 8F.28") of the MST deck:
# the heading text box (shape 3, "TextBox 3") reads "color components" again
# even though this is actually the *second* time through the color/gray-edge
# loop (the first coloring happens on slides 24-25), so it should read
# "re-color components" instead. The text box uses shape-autofit sizing, so
# its width (ext cx) needs to grow to fit the now-longer heading as well.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(28)
$sh = $s.Shapes.Item(3)

# Widen the auto-fit text box to accommodate the longer heading
# (4725998 EMU -> 5541250 EMU, i.e. 372.1258pt -> 436.3189pt; height unchanged).
$sh.Width = 436.3188976377953

$tr = $sh.TextFrame.TextRange

# "color components" -> "re-color components"
# Remove the leading "color " (6 chars) and insert "re-color " in front of
# the remaining "components", mirroring how PowerPoint splits runs when text
# is edited in place (both new runs keep the original 44pt / +mj-lt formatting).
$tr.Characters(1, 6).Delete() | Out-Null
$tr.InsertBefore("re-color ") | Out-Null
